# This script removes the "MDL" (Mandalay, Myanmar) row from the Colos
# data sheet. The row was removed entirely from the source data set, so
# all subsequent rows shift up by one and the used range shrinks by one
# row (from A1:H338 to A1:H337).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Locate the row whose colo code (column A) is "MDL" and delete it,
# shifting the rows below it upward.
$found = $ws.Cells.Find("MDL", [System.Reflection.Missing]::Value, -4163, 1)

if ($found -ne $null) {
    $rowNum = $found.Row
    $ws.Rows.Item($rowNum).Delete()
} else {
    # Fallback: we know from the data layout that the MDL row is row 222.
    $ws.Rows.Item(222).Delete()
}
